$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to be treated as text so values like
# "1.002" or "81.50" are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.207.36'
$ws.Range("E2").Value = '  -3.48%  '
$ws.Range("D3").Value = '1.807.32'
$ws.Range("E3").Value = '  -3.63%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '310.49'
$ws.Range("E5").Value = '  -1.67%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.4209'
$ws.Range("E7").Value = '  -2.21%  '
$ws.Range("D8").Value = '0.3552'
$ws.Range("E8").Value = '  -3.48%  '
$ws.Range("D9").Value = '0.07115'
$ws.Range("E9").Value = '  -3.93%  '
$ws.Range("D10").Value = '0.8489'
$ws.Range("E10").Value = '  -3.45%  '
$ws.Range("D11").Value = '20.17'
$ws.Range("E11").Value = '  -4.31%  '
$ws.Range("D12").Value = '1.776.50'
$ws.Range("E12").Value = '  -7.28%  '
$ws.Range("D13").Value = '5.317'
$ws.Range("E13").Value = '  -2.89%  '
$ws.Range("D14").Value = '6.373'
$ws.Range("E14").Value = '  -3.59%  '
$ws.Range("D15").Value = '0.06882'
$ws.Range("E15").Value = '  -1.19%  '
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = '81.50'
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").Value = '0.000008771'
$ws.Range("E18").Value = '  -3.90%  '
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("D20").Value = '15.08'
$ws.Range("E20").Value = '  -3.07%  '
$ws.Range("D21").Value = '27.328.08'
$ws.Range("E21").Value = '  -3.39%  '
$ws.Range("D22").Value = '5.086'
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("E23").Value = '  -0.75%  '
$ws.Range("D24").Value = '2.037.39'
$ws.Range("E24").Value = '  -6.10%  '
$ws.Range("D25").Value = '1.964'
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("D26").Value = '153.64'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '18.23'
$ws.Range("E27").Value = '  -2.52%  '
$ws.Range("D28").Value = '5.053'
$ws.Range("E28").Value = '  -6.33%  '
$ws.Range("D29").Value = '113.20'
$ws.Range("E29").Value = '  -3.60%  '
$ws.Range("D30").Value = '1.709'
$ws.Range("E30").Value = '  -8.50%  '
$ws.Range("D31").Value = '0.08888'
$ws.Range("E31").Value = '  -0.70%  '
$ws.Range("D32").Value = '0.7414'
$ws.Range("E32").Value = '  -5.83%  '
$ws.Range("D33").Value = '4.458'
$ws.Range("E33").Value = '  -5.19%  '
$ws.Range("D34").Value = '2.911'
$ws.Range("E34").Value = '  -2.00%  '
$ws.Range("D35").Value = '1.102'
$ws.Range("E35").Value = '  -6.25%  '
$ws.Range("D36").Value = '1.002'
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").Value = '1.068'
$ws.Range("E37").Value = '  -5.49%  '
$ws.Range("D38").Value = '0.05204'
$ws.Range("E38").Value = '  -4.71%  '
$ws.Range("D39").Value = '0.01897'
$ws.Range("E39").Value = '  -3.14%  '
$ws.Range("D40").Value = '0.1636'
$ws.Range("E40").Value = '  -3.18%  '
$ws.Range("D41").Value = '2.709'
$ws.Range("E41").Value = '  -6.34%  '
$ws.Range("D42").Value = '0.4967'
$ws.Range("E42").Value = '  -3.66%  '
$ws.Range("D43").Value = '6.289'
$ws.Range("E43").Value = '  -8.33%  '
$ws.Range("D44").Value = '8.177'
$ws.Range("E44").Value = '  -4.48%  '
$ws.Range("D45").Value = '105.04'
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("D46").Value = '10.18'
$ws.Range("E46").Value = '  -3.50%  '
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("D48").Value = '0.06384'
$ws.Range("E48").Value = '  -3.04%  '
$ws.Range("D49").Value = '0.4556'
$ws.Range("E49").Value = '  -4.16%  '
$ws.Range("D50").Value = '1.594'
$ws.Range("E50").Value = '  -3.54%  '
$ws.Range("D51").Value = '62.84'
$ws.Range("E51").Value = '  -3.94%  '

# Restore default (General) cell style so no style index is left behind
# on the edited cells, matching the original workbook formatting.
$ws.Range("D2:E51").Style = "Normal"
